$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The review for row 19 (innaplutov1@gmail.com / "App is working as expected after
# update..") is being removed; everything below it shifts up by one row. Capture
# the hyperlinks that live on the rows which will shift (old C21 -> C20,
# old C22 -> C21) so we can recreate them in their new spot, since a plain row
# delete does not relocate existing Hyperlink objects in this engine.
$hls = $ws.Hyperlinks
$hlList = @()
foreach ($hl in $hls) {
    $hlList += $hl
}

$hlC21 = $hlList[19]
$hlC22 = $hlList[20]
$addrC21 = $hlC21.Address
$dispC21 = $hlC21.TextToDisplay
$addrC22 = $hlC22.Address
$dispC22 = $hlC22.TextToDisplay

# Remove the two hyperlink objects that would otherwise be left pointing at the
# wrong (stale) cells once the row above them is deleted.
$hlC22.Delete()
$hlC21.Delete()

# Delete the whole row 19 - remaining rows (20-24) shift up to (19-23).
$ws.Rows(19).Delete()

# Recreate the hyperlinks at their new, shifted-up locations.
$ws.Hyperlinks.Add($ws.Range("C20"), $addrC21, "", "", $dispC21)
$ws.Hyperlinks.Add($ws.Range("C21"), $addrC22, "", "", $dispC22)

# Match the author's final selection/view state.
$ws.Range("A19").Select()
